$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 361 (shifts existing rows 361:388 down to 362:389)
$ws.Rows(361).Insert()

# Populate the new row 361 with a new data record (same template as the
# surrounding rows, with its own Date/Volumen/Precio values)
$ws.Cells.Item(361, 1).Value = 3
$ws.Cells.Item(361, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(361, 3).Value = "Coquimbo"
$ws.Cells.Item(361, 4).Value = 44783
$ws.Cells.Item(361, 5).Value = 5
$ws.Cells.Item(361, 6).Value = 100112009
$ws.Cells.Item(361, 7).Value = "Acelga"
$ws.Cells.Item(361, 8).Value = "Sin especificar"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 230
$ws.Cells.Item(361, 11).Value = 3000
$ws.Cells.Item(361, 12).Value = 3300
$ws.Cells.Item(361, 13).Value = 3157
$ws.Cells.Item(361, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(361, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(361, 16).Value = 526
$ws.Cells.Item(361, 17).Value = 6
$ws.Cells.Item(361, 18).Value = "Hortaliza"
